$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Hoja1)
$ws.Name = "Hoja1"

# --- Header row updates ---
# Var1 -> Account
$ws.Range("C1").Value = "Account"
# New "AccountKO" header column
$ws.Range("F1").Value = "AccountKO"

# --- Existing data updates ---
# Row 2: Run flag Y -> N (api test now expected to fail on this first approach)
$ws.Range("B2").Value = "N"

# --- New row 5: sample_apitests_excel test case ---
$ws.Range("A5").Value = "sample_apitests_excel"
$ws.Range("B5").Value = "Y"
$ws.Range("C5").Value = "anxoportela"
$ws.Range("D5").Value = "b"
$ws.Range("E5").Value = "c"
$ws.Range("F5").Value = "l324dsg34hguisito"

# --- Apply the Arial 11 font to all populated cells (matches the rest of the workbook look) ---
$dataRanges = @("A1:D1", "E1:F1", "A2:D2", "E2:E2", "A3:B3", "A4:B4", "A5:D5", "E5:F5")
foreach ($r in $dataRanges) {
    $ws.Range($r).Font.Name = "Arial"
    $ws.Range($r).Font.Size = 11
}

# --- Row heights: every row now lines up at 13.8pt ---
$ws.Range("A1:F1").RowHeight = 13.8
$ws.Range("A2:F2").RowHeight = 13.8
$ws.Range("A3:F3").RowHeight = 13.8
$ws.Range("A4:F4").RowHeight = 13.8
$ws.Range("A5:F5").RowHeight = 13.8

# --- Selection moves to F6, just past the new data ---
$ws.Range("F6").Select() | Out-Null
